$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5490601.8474
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = 0.4637
$ws.Range("D3").Value = 0.9268

$ws.Range("C4").Value = 3.2394
$ws.Range("D4").Value = 0.3562

$ws.Range("C5").Value = 1.0849
$ws.Range("D5").Value = 0.7806999999999999

$ws.Range("C6").Value = 13.0433
$ws.Range("D6").Value = 0.0045

$ws.Range("C7").Value = 386.978
$ws.Range("D7").Value = 0
